$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Update capital-cost numbers in column C, rows 3-17 (re-evaluated
#    results from the author's updated model run).
# -----------------------------------------------------------------
$ws.Range("C3").Value  = 405.8710005627602
$ws.Range("C4").Value  = 62.47830084552534
$ws.Range("C5").Value  = 16.23484002251041
$ws.Range("C6").Value  = 36.52839005064841
$ws.Range("C7").Value  = 18.26419502532421
$ws.Range("C8").Value  = 539.3767265067685
$ws.Range("C9").Value  = 53.93767265067685
$ws.Range("C10").Value = 53.93767265067685
$ws.Range("C11").Value = 107.8753453013537
$ws.Range("C12").Value = 161.8130179520305
$ws.Range("C13").Value = 53.93767265067685
$ws.Range("C14").Value = 431.5013812054149
$ws.Range("C15").Value = 970.8781077121834
$ws.Range("C16").Value = 48.54390538560917
$ws.Range("C17").Value = 1019.422013097793

# -----------------------------------------------------------------
# 2) Raw-materials / by-products table (rows 21-28): the two new rows
#    "Tridecane" and "CSL" now sit at the top (right after the row 20
#    header), pushing DAP/Glucose/Salt/Wastewater/Process water down
#    so the whole block shifts. Unmerge the old A-column groupings
#    before rewriting, then re-merge to match the new grouping.
# -----------------------------------------------------------------
$ws.Range("A21:A23").UnMerge()
$ws.Range("A25:A28").UnMerge()

$rows = @(
    @{ Row = 21; A = "Raw materials";             B = "Tridecane";     C = 878.1550799999999;  D = 0.02944919024573181 },
    @{ Row = 22; A = $null;                       B = "CSL";           C = 51.528108;           D = 6.106799307209393 },
    @{ Row = 23; A = $null;                       B = "DAP";           C = 895.3915949999999;   D = 14.44643600469633 },
    @{ Row = 24; A = $null;                       B = "Glucose";       C = 240.404025;          D = 231.5793396277901 },
    @{ Row = 25; A = $null;                       B = "Salt";          C = 136.07775;           D = 51.8569443828605 },
    @{ Row = 26; A = "By-products and credits";   B = "Wastewater";    C = -1.122754726231208;  D = -43.07057236753541 },
    @{ Row = 27; A = "Raw materials";              B = "Process water"; C = 0.320236305;         D = 12.20091644845022 },
    @{ Row = 28; A = $null;                       B = "Natural gas";   C = 197.76633;           D = 12.40122128402908 }
)

foreach ($r in $rows) {
    $row = $r.Row
    if ($null -ne $r.A) {
        $ws.Range("A$row").Value = $r.A
    } else {
        $ws.Range("A$row").Value = $null
    }
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
}

$ws.Range("A21:A25").Merge()
$ws.Range("A27:A28").Merge()

# Total variable operating cost (row 29, column D)
$ws.Range("D29").Value = 371.6916786128169

# -----------------------------------------------------------------
# 3) Labor / overhead table numbers (rows 35-36).
# -----------------------------------------------------------------
$ws.Range("C35").Value = 12.17613001688281
$ws.Range("D35").Value = 11.68908481620749
$ws.Range("C36").Value = 2.841097003939321
$ws.Range("D36").Value = 2.727453123781749
